$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 5000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -4685
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 5000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -3908
$ws.Range("N79").ClearContents()

$ws.Range("H110").Value = 35996.5
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 35996.5
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 35996.5
$ws.Range("N110").Value = -44176.5

$ws.Range("H111").Value = 7697.8237
$ws.Range("I111").Value = 9798.362999999999
$ws.Range("J111").Value = 3846.8333
$ws.Range("K111").Value = 29395.089
$ws.Range("L111").Value = 11540.4999
$ws.Range("M111").Value = -26328.089
$ws.Range("N111").Value = -17674.4999

$ws.Range("H112").Value = 30041.361
$ws.Range("I112").Value = 1015.5
$ws.Range("J112").Value = 35846.535
$ws.Range("K112").Value = 3046.5
$ws.Range("L112").Value = 107539.605
$ws.Range("M112").Value = -1938.5
$ws.Range("N112").Value = -109755.605

$ws.Range("H137").Value = 1522.1666
$ws.Range("I137").Value = 1404.6552
$ws.Range("J137").Value = 1701.5264
$ws.Range("K137").Value = 4213.9656
$ws.Range("L137").Value = 5104.5792
$ws.Range("M137").Value = -1663.9656
$ws.Range("N137").Value = -10204.5792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1558.4286
$ws.Range("I2").Value = 1293.6923
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 1293.6923
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -1180.6923
$ws.Range("N2").Value = -5226

$ws.Range("H12").Value = 209.42857
$ws.Range("I12").Value = 44.333332
$ws.Range("J12").Value = 1200
$ws.Range("K12").Value = 44.333332
$ws.Range("L12").Value = 1200
$ws.Range("M12").Value = 128.666668
$ws.Range("N12").Value = -1546

$ws.Range("H116").Value = 1558.4286
$ws.Range("I116").Value = 1293.6923
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 1293.6923
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 1000.3077
$ws.Range("N116").Value = -9588

$ws.Range("H132").Value = 3040.08
$ws.Range("I132").Value = 2731.85
$ws.Range("J132").Value = 4273
$ws.Range("K132").Value = 8195.549999999999
$ws.Range("L132").Value = 12819
$ws.Range("M132").Value = -5665.549999999999
$ws.Range("N132").Value = -17879

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1558.4286
$ws.Range("I3").Value = 1293.6923
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 1293.6923
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -1179.6923
$ws.Range("N3").Value = -5228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2001
$ws.Range("I3").Value = 2500
$ws.Range("J3").Value = 1751.5
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 1751.5
$ws.Range("M3").Value = -2387
$ws.Range("N3").Value = -1977.5

$ws.Range("H22").Value = 230
$ws.Range("I22").Value = 230
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 230
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 120

$ws.Range("H93").Value = 11498.5
$ws.Range("I93").Value = 11498.5
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 11498.5
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -9626.5

$ws.Range("H103").Value = 15797.2
$ws.Range("I103").Value = 15797.2
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 15797.2
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -14625.2

$ws.Range("H133").Value = 49500
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 49500
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 49500
$ws.Range("N133").Value = -54560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 119.77778
$ws.Range("I2").Value = 153.33333
$ws.Range("J2").Value = 52.666668
$ws.Range("K2").Value = 919.9999799999999
$ws.Range("L2").Value = 316.000008
$ws.Range("M2").Value = -806.9999799999999
$ws.Range("N2").Value = -542.000008

$ws.Range("H23").Value = 314.44
$ws.Range("I23").Value = 325.8
$ws.Range("J23").Value = 306.86667
$ws.Range("K23").Value = 977.4000000000001
$ws.Range("L23").Value = 920.60001
$ws.Range("M23").Value = -742.4000000000001
$ws.Range("N23").Value = -1390.60001

$ws.Range("H25").Value = 835.25
$ws.Range("I25").Value = 671.25
$ws.Range("J25").Value = 999.25
$ws.Range("K25").Value = 2013.75
$ws.Range("L25").Value = 2997.75
$ws.Range("M25").Value = -1844.75
$ws.Range("N25").Value = -3335.75

$ws.Range("H30").Value = 835.25
$ws.Range("I30").Value = 671.25
$ws.Range("J30").Value = 999.25
$ws.Range("K30").Value = 2013.75
$ws.Range("L30").Value = 2997.75
$ws.Range("M30").Value = -1911.75
$ws.Range("N30").Value = -3201.75

$ws.Range("H32").Value = 3654.4285
$ws.Range("I32").Value = 4500
$ws.Range("J32").Value = 3513.5
$ws.Range("K32").Value = 13500
$ws.Range("L32").Value = 10540.5
$ws.Range("M32").Value = -13217
$ws.Range("N32").Value = -11106.5

$ws.Range("H38").Value = 36.75
$ws.Range("I38").Value = 33.444443
$ws.Range("J38").Value = 41
$ws.Range("K38").Value = 100.333329
$ws.Range("L38").Value = 123
$ws.Range("M38").Value = 246.666671
$ws.Range("N38").Value = -817

$ws.Range("H75").Value = 1183.6
$ws.Range("I75").Value = 1301
$ws.Range("J75").Value = 1007.5
$ws.Range("K75").Value = 3903
$ws.Range("L75").Value = 3022.5
$ws.Range("M75").Value = -2905
$ws.Range("N75").Value = -5018.5

$ws.Range("H78").Value = 1183.6
$ws.Range("I78").Value = 1301
$ws.Range("J78").Value = 1007.5
$ws.Range("K78").Value = 11709
$ws.Range("L78").Value = 9067.5
$ws.Range("M78").Value = -6717
$ws.Range("N78").Value = -19051.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 29000
$ws.Range("I46").Value = 8000
$ws.Range("J46").Value = 50000
$ws.Range("K46").Value = 8000
$ws.Range("L46").Value = 50000
$ws.Range("M46").Value = -7844
$ws.Range("N46").Value = -50312

$ws.Range("H122").Value = 2368.9736
$ws.Range("I122").Value = 2207.5938
$ws.Range("J122").Value = 3229.6667
$ws.Range("K122").Value = 6622.7814
$ws.Range("L122").Value = 9689.000100000001
$ws.Range("M122").Value = -4172.7814
$ws.Range("N122").Value = -14589.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 131.04546
$ws.Range("I55").Value = 126.23529
$ws.Range("J55").Value = 147.4
$ws.Range("K55").Value = 126.23529
$ws.Range("L55").Value = 147.4
$ws.Range("M55").Value = 46.76470999999999
$ws.Range("N55").Value = -493.4

$ws.Range("H110").Value = 45999.75
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 45999.75
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 45999.75
$ws.Range("N110").Value = -54179.75

$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 38135.934
$ws.Range("I69").Value = 34814.5
$ws.Range("J69").Value = 51421.668
$ws.Range("K69").Value = 34814.5
$ws.Range("L69").Value = 51421.668
$ws.Range("M69").Value = -34065.5
$ws.Range("N69").Value = -52919.668

$ws.Range("H72").Value = 38135.934
$ws.Range("I72").Value = 34814.5
$ws.Range("J72").Value = 51421.668
$ws.Range("K72").Value = 104443.5
$ws.Range("L72").Value = 154265.004
$ws.Range("M72").Value = -100699.5
$ws.Range("N72").Value = -161753.004

$ws.Range("H81").Value = 10749.4375
$ws.Range("I81").Value = 51799.5
$ws.Range("J81").Value = 4885.143
$ws.Range("K81").Value = 103599
$ws.Range("L81").Value = 9770.286
$ws.Range("M81").Value = -102538
$ws.Range("N81").Value = -11892.286

$ws.Range("H84").Value = 10749.4375
$ws.Range("I84").Value = 51799.5
$ws.Range("J84").Value = 4885.143
$ws.Range("K84").Value = 517995
$ws.Range("L84").Value = 48851.43
$ws.Range("M84").Value = -512691
$ws.Range("N84").Value = -59459.43

$ws.Range("H136").Value = 2323.5
$ws.Range("I136").Value = 2309.2666
$ws.Range("J136").Value = 2376.875
$ws.Range("K136").Value = 6927.7998
$ws.Range("L136").Value = 7130.625
$ws.Range("M136").Value = -4377.7998
$ws.Range("N136").Value = -12230.625
